$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d17 = @'
Background
id="Par1">An outbreak of coronavirus disease 2019 (COVID-19) caused by SARS-CoV-2 was first detected in Wuhan, Hubei, China.

 People of all ages are susceptible to SARS-CoV-2 infection.

 No information on severe pediatric patients with COVID-19 has been reported.

 We aimed to describe the clinical features of severe pediatric patients with COVID-19.
Methods
id="Par2">We included eight severe or critically ill patients with COVID-19 who were treated at the Intensive Care Unit (ICU), Wuhan Children’s Hospital from January 24 to February 24. We collected information including demographic data, symptoms, imaging data, laboratory findings, treatments and clinical outcomes of the patients with severe COVID-19.
Results
id="Par3">The onset age of the eight patients ranged from 2 months to 15 years; six were boys.

 The most common symptoms were polypnea (8/8), followed by fever (6/8) and cough (6/8).

 Chest imaging showed multiple patch-like shadows in seven patients and ground-glass opacity in six.

 Laboratory findings revealed normal or increased whole blood counts (7/8), increased C-reactive protein, procalcitonin and lactate dehydrogenase (6/8), and abnormal liver function (4/8).

 Other findings included decreased CD16 + CD56 (4/8) and Th/Ts*(1/8), increased CD3 (2/8), CD4 (4/8) and CD8 (1/8), IL-6 (2/8), IL-10 (5/8) and IFN-γ (2/8).

 Treatment modalities were focused on symptomatic and respiratory support.

 Two critically ill patients underwent invasive mechanical ventilation.

 Up to February 24, 2020, three patients remained under treatment in ICU, the other five recovered and were discharged home.


Conclusions
id="Par4">In this series of severe pediatric patients in Wuhan, polypnea was the most common symptom, followed by fever and cough.

 Common imaging changes included multiple patch-like shadows and ground-glass opacity; and a cytokine storm was found in these patients, which appeared more serious in critically ill patients.



'@

$e17 = @'
[Dan%Sun%NULL%0,  Hui%Li%NULL%0,  Xiao-Xia%Lu%NULL%0,  Han%Xiao%NULL%0,  Jie%Ren%NULL%0,  Fu-Rong%Zhang%792523496@qq.com%0,  Zhi-Sheng%Liu%liuzsc@126.com%0]
'@

$ws.Range("D17").Value2 = $d17
$ws.Range("E17").Value2 = $e17
